$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (rows 336-337),
# pushing the existing rows (336:370) down to (338:372).
$ws.Rows("336:337").Insert()

# New week's data (Primera) - row 336
$ws.Range("A336").Value = 8
$ws.Range("B336").Value = "Terminal La Palmera de La Serena"
$ws.Range("C336").Value = "Coquimbo"
$ws.Range("D336").Value = 44449
$ws.Range("E336").Value = 4
$ws.Range("F336").Value = 100112008
$ws.Range("G336").Value = "Coliflor"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 3600
$ws.Range("K336").Value = 650
$ws.Range("L336").Value = 700
$ws.Range("M336").Value = 675
$ws.Range("N336").Value = "`$/unidad"
$ws.Range("O336").Value = "Provincia del Elquí"
$ws.Range("P336").Value = 675
$ws.Range("Q336").Value = 1
$ws.Range("R336").Value = "Hortaliza"

# New week's data (Segunda) - row 337
$ws.Range("A337").Value = 8
$ws.Range("B337").Value = "Terminal La Palmera de La Serena"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44449
$ws.Range("E337").Value = 4
$ws.Range("F337").Value = 100112008
$ws.Range("G337").Value = "Coliflor"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "Segunda"
$ws.Range("J337").Value = 1900
$ws.Range("K337").Value = 550
$ws.Range("L337").Value = 600
$ws.Range("M337").Value = 575
$ws.Range("N337").Value = "`$/unidad"
$ws.Range("O337").Value = "Provincia del Elquí"
$ws.Range("P337").Value = 575
$ws.Range("Q337").Value = 1
$ws.Range("R337").Value = "Hortaliza"
